$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(0.926183767058954, 0.99856345654456, 0.998673959887287, 0.938615393115642, 0.00044201337090447, 0.000552516713630587, 0.989170672412841, 0.267307586054478, 0.000386761699541411, 0.948118680590088, 0.000663020056356705, 0.999226476600917, 0.952759820984585, 0.000386761699541411, 0.999723741643185, 0.242610088955191, 0.999336979943643, 0.997624178131388, 0.0312724459914912, 0.0125421293994143, 0.013205149455771, 0.00154704679816564, 0.000331510028178352)
$row3 = @(0.0079009890049174, 0.0000552516713630587, 0.0000552516713630587, 0.000607768384993646, 0.000331510028178352, 0.00237582186861153, 0.00044201337090447, 0.000165755014089176, 0.00088402674180894, 0.000276258356815294, 0.000110503342726117, 0.0000552516713630587, 0, 0, 0, 0.000939278413171998, 0, 0, 0, 0.000276258356815294, 0.985358307088789, 0.0686225758329189, 0.00104978175589812)
$row4 = @(0.0518813194099121, 0.000663020056356705, 0.000386761699541411, 0.0588982816730206, 0.998397701530471, 0.000718271727719763, 0.00911652577490469, 0.730482347090999, 0, 0.0502790209403834, 0.998176694845019, 0.00044201337090447, 0.0453616221890712, 0.999502734957732, 0.000221006685452235, 0.748991656997624, 0.000497265042267529, 0.00193380849770706, 0.96640698381126, 0.955467152881375, 0.00110503342726117, 0.00243107353997458, 0.000331510028178352)
$row5 = @(0.0123763743853252, 0.000718271727719763, 0.00088402674180894, 0.000221006685452235, 0.000828775070445881, 0.996242886347312, 0.000221006685452235, 0, 0.998673959887287, 0.000331510028178352, 0.000994530084535057, 0.000276258356815294, 0.000165755014089176, 0.0000552516713630587, 0, 0.000221006685452235, 0.000165755014089176, 0.000221006685452235, 0.0000552516713630587, 0.0315487043483065, 0.000221006685452235, 0.925023481960329, 0.998287198187745)

$rows = @(2, 3, 4, 5)
$allRows = @($row2, $row3, $row4, $row5)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    $vals = $allRows[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 2).Value = $vals[$c]
    }
}
